$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps being stored as text (it holds values like
# "37.340.39" that are not valid numbers, so format the whole column as text
# before writing any of the new values).
$ws.Range("D2:D51").NumberFormat = "@"

# Rows 42 and 44 swap their Coin/Link data and get new Price/Volume values
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.486.40"
$ws.Range("E42").Value = "  +2.85%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0951"
$ws.Range("E44").Value = "  -3.02%  "

$ws.Range("D2").Value = '37.340.39'
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").Value = '2.050.94'
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '230.39'
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("D6").Value = '0.621'
$ws.Range("E6").Value = '  -0.70%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '57.06'
$ws.Range("E8").Value = '  -4.01%  '
$ws.Range("E9").Value = '  -2.47%  '
$ws.Range("D10").Value = '0.0769'
$ws.Range("E10").Value = '  -2.64%  '
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").Value = '14.71'
$ws.Range("E12").Value = '  -0.95%  '
$ws.Range("D13").Value = '2.354.23'
$ws.Range("E13").Value = '  -1.30%  '
$ws.Range("D14").Value = '20.55'
$ws.Range("E14").Value = '  -3.23%  '
$ws.Range("D15").Value = '0.756'
$ws.Range("E15").Value = '  -2.67%  '
$ws.Range("D16").Value = '5.24'
$ws.Range("E16").Value = '  -2.41%  '
$ws.Range("D17").Value = '2.047.13'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '37.312.15'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("E19").Value = '  -2.89%  '
$ws.Range("D20").Value = '69.80'
$ws.Range("E20").Value = '  -2.55%  '
$ws.Range("D21").Value = '0.0₃0824'
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("D22").Value = '226.62'
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("E25").Value = '  -3.90%  '
$ws.Range("D26").Value = '9.48'
$ws.Range("E26").Value = '  +3.07%  '
$ws.Range("D27").Value = '168.76'
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("E28").Value = '  -3.61%  '
$ws.Range("D29").Value = '19.15'
$ws.Range("E29").Value = '  -1.89%  '
$ws.Range("D30").Value = '1.35'
$ws.Range("E30").Value = '  -5.15%  '
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").Value = '4.53'
$ws.Range("E32").Value = '  -4.23%  '
$ws.Range("D33").Value = '0.0624'
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("D34").Value = '4.56'
$ws.Range("E34").Value = '  -4.10%  '
$ws.Range("D35").Value = '2.49'
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").Value = '3.28'
$ws.Range("E37").Value = '  -3.91%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("D39").Value = '5.26'
$ws.Range("E39").Value = '  -2.54%  '
$ws.Range("E40").Value = '  +4.45%  '
$ws.Range("D41").Value = '98.03'
$ws.Range("E41").Value = '  -1.46%  '
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("E45").Value = '  +2.89%  '
$ws.Range("D46").Value = '16.59'
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("D47").Value = '4.03'
$ws.Range("E47").Value = '  -2.91%  '
$ws.Range("E48").Value = '  -3.46%  '
$ws.Range("D49").Value = '7.23'
$ws.Range("E49").Value = '  -2.37%  '
$ws.Range("E50").Value = '  -3.16%  '
$ws.Range("D51").Value = '2.240.45'
$ws.Range("E51").Value = '  -1.30%  '